$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 27 de Marzo de 2020 a las 14:12"

# Row 4
$ws.Cells.Item(4, 2).Value = 85755
$ws.Cells.Item(4, 3).Value = 320
$ws.Cells.Item(4, 5).Value = 82583

# Row 11
$ws.Cells.Item(11, 2).Value = 12311
$ws.Cells.Item(11, 3).Value = 500
$ws.Cells.Item(11, 5).Value = 11207
$ws.Cells.Item(11, 7).Value = 15
$ws.Cells.Item(11, 8).Value = 207

# Row 14
$ws.Cells.Item(14, 2).Value = 8603
$ws.Cells.Item(14, 3).Value = 1172
$ws.Cells.Item(14, 5).Value = 8054
$ws.Cells.Item(14, 7).Value = 112
$ws.Cells.Item(14, 8).Value = 546

# Row 15
$ws.Cells.Item(15, 2).Value = 7393
$ws.Cells.Item(15, 3).Value = 484
$ws.Cells.Item(15, 5).Value = 7110

# Row 19
$ws.Cells.Item(19, 2).Value = 3677
$ws.Cells.Item(19, 3).Value = 305
$ws.Cells.Item(19, 5).Value = 3655

# Row 22
$ws.Cells.Item(22, 1).Value = "Suecia"
$ws.Cells.Item(22, 2).Value = 3046
$ws.Cells.Item(22, 3).Value = 206
$ws.Cells.Item(22, 4).Value = 16
$ws.Cells.Item(22, 5).Value = 2938
$ws.Cells.Item(22, 6).Value = 209
$ws.Cells.Item(22, 7).Value = 15
$ws.Cells.Item(22, 8).Value = 92

# Row 23
$ws.Cells.Item(23, 1).Value = "Israel"
$ws.Cells.Item(23, 2).Value = 3035
$ws.Cells.Item(23, 3).Value = 342
$ws.Cells.Item(23, 4).Value = 79
$ws.Cells.Item(23, 5).Value = 2945
$ws.Cells.Item(23, 6).Value = 49
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = 11

# Row 24
$ws.Cells.Item(24, 1).Value = "Brasil"
$ws.Cells.Item(24, 2).Value = 2988
$ws.Cells.Item(24, 3).Value = 3
$ws.Cells.Item(24, 4).Value = 6
$ws.Cells.Item(24, 5).Value = 2905
$ws.Cells.Item(24, 6).Value = 296

# Row 34
$ws.Cells.Item(34, 1).Value = "Polonia"
$ws.Cells.Item(34, 2).Value = 1289
$ws.Cells.Item(34, 3).Value = 68
$ws.Cells.Item(34, 4).Value = 7
$ws.Cells.Item(34, 5).Value = 1266
$ws.Cells.Item(34, 6).Value = 3
$ws.Cells.Item(34, 8).Value = 16

# Row 35
$ws.Cells.Item(35, 1).Value = "Pakistan"
$ws.Cells.Item(35, 2).Value = 1252
$ws.Cells.Item(35, 3).Value = 51
$ws.Cells.Item(35, 4).Value = 23
$ws.Cells.Item(35, 5).Value = 1220
$ws.Cells.Item(35, 6).Value = 7
$ws.Cells.Item(35, 8).Value = 9

# Row 37
$ws.Cells.Item(37, 1).Value = "Arabia Saudita"
$ws.Cells.Item(37, 2).Value = 1104
$ws.Cells.Item(37, 3).Value = 92
$ws.Cells.Item(37, 4).Value = 35
$ws.Cells.Item(37, 5).Value = 1066
$ws.Cells.Item(37, 6).Value = 6
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 3

# Row 38
$ws.Cells.Item(38, 1).Value = "Indonesia"
$ws.Cells.Item(38, 2).Value = 1046
$ws.Cells.Item(38, 3).Value = 153
$ws.Cells.Item(38, 4).Value = 46
$ws.Cells.Item(38, 5).Value = 913
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 9
$ws.Cells.Item(38, 8).Value = 87

# Row 39
$ws.Cells.Item(39, 1).Value = "Finlandia"
$ws.Cells.Item(39, 2).Value = 1038
$ws.Cells.Item(39, 3).Value = 80
$ws.Cells.Item(39, 4).Value = 10
$ws.Cells.Item(39, 5).Value = 1023
$ws.Cells.Item(39, 6).Value = 24
$ws.Cells.Item(39, 8).Value = 5

# Row 40
$ws.Cells.Item(40, 1).Value = "Rusia"
$ws.Cells.Item(40, 2).Value = 1036
$ws.Cells.Item(40, 3).Value = 196
$ws.Cells.Item(40, 4).Value = 45
$ws.Cells.Item(40, 5).Value = 988
$ws.Cells.Item(40, 6).Value = 8

# Row 43
$ws.Cells.Item(43, 1).Value = "Islandia"
$ws.Cells.Item(43, 2).Value = 890
$ws.Cells.Item(43, 3).Value = 88
$ws.Cells.Item(43, 4).Value = 82
$ws.Cells.Item(43, 5).Value = 806
$ws.Cells.Item(43, 6).Value = 17
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 2

# Row 44
$ws.Cells.Item(44, 2).Value = 843
$ws.Cells.Item(44, 3).Value = 116
$ws.Cells.Item(44, 5).Value = 750

# Row 45
$ws.Cells.Item(45, 1).Value = "Filipinas"
$ws.Cells.Item(45, 2).Value = 803
$ws.Cells.Item(45, 3).Value = 96
$ws.Cells.Item(45, 4).Value = 31
$ws.Cells.Item(45, 6).Value = 1
$ws.Cells.Item(45, 7).Value = 9
$ws.Cells.Item(45, 8).Value = 54

# Row 50
$ws.Cells.Item(50, 5).Value = 504
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 13

# Row 62
$ws.Cells.Item(62, 6).Value = 25

# Row 94
$ws.Cells.Item(94, 1).Value = "Kazajistan"
$ws.Cells.Item(94, 2).Value = 135
$ws.Cells.Item(94, 4).Value = 3
$ws.Cells.Item(94, 5).Value = 131
$ws.Cells.Item(94, 8).Value = 1

# Row 95
$ws.Cells.Item(95, 1).Value = "Oman"
$ws.Cells.Item(95, 2).Value = 131
$ws.Cells.Item(95, 3).Value = 22
$ws.Cells.Item(95, 4).Value = 23
$ws.Cells.Item(95, 5).Value = 108
$ws.Cells.Item(95, 8).Value = 0

# Row 127
$ws.Cells.Item(127, 1).Value = "Macao"
$ws.Cells.Item(127, 2).Value = 34
$ws.Cells.Item(127, 3).Value = 1
$ws.Cells.Item(127, 4).Value = 10
$ws.Cells.Item(127, 5).Value = 24

# Row 128
$ws.Cells.Item(128, 1).Value = "Monaco"
$ws.Cells.Item(128, 4).Value = 1
$ws.Cells.Item(128, 5).Value = 32
